$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21 - this pushes the existing
# "frontPanelThickness" row (and everything below it) down by one,
# turning old row 21 into row 22.
$ws.Rows("21:21").Insert() | Out-Null

# Populate the newly inserted row 21 with the new parameter:
# fresnelGrooveWidth = 0.07 in (70 thou)
$ws.Range("A21").Value = "fresnelGrooveWidth"
$ws.Range("B21").Value = 0.07
$ws.Range("C21").Value = "in"

# Match the number format used by similar "thickness"-style parameters
# (style index 1 / numFmtId 164, "0.000") rather than the format that
# got copied down from the row above during the insert.
$ws.Range("B21").NumberFormat = "0.000"

# frontPanelThickness (now on row 22) is derived by subtracting the new
# groove width from 2, instead of subtracting the fixed 1/16 literal.
$ws.Range("B8").Formula = "=2-B21"

# Update the active selection to match the authored workbook.
$ws.Range("D15").Select() | Out-Null
